$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I0 (col I) and IF (col J), copying style from existing header cell H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row data for new columns I (I0) and J (IF)
$data = @{
  2 = @(7,7)
  3 = @(8,8)
  4 = @(7,7)
  5 = @(6,6)
  6 = @(7,7)
  7 = @(7,7)
  8 = @(11,11)
  9 = @(7,7)
  10 = @(6,6)
  11 = @(8,8)
  12 = @(5,6)
  13 = @(8,8)
  14 = @(8,8)
  15 = @(7,7)
  16 = @(8,8)
  17 = @(11,11)
  18 = @(7,7)
  19 = @(5,6)
  20 = @(7,7)
  21 = @(8,8)
  22 = @(8,8)
  23 = @(7,7)
  24 = @(6,6)
  25 = @(7,7)
  26 = @(7,7)
  27 = @(6,6)
  28 = @(8,8)
  29 = @(6,7)
  30 = @(7,7)
  31 = @(9,9)
  32 = @(7,7)
  33 = @(4,4)
  34 = @(8,8)
  35 = @(7,7)
  36 = @(6,6)
  37 = @(10,10)
  38 = @(5,5)
  39 = @(6,6)
  40 = @(10,11)
  41 = @(6,6)
  42 = @(6,7)
  43 = @(7,7)
  44 = @(7,7)
  45 = @(6,6)
  46 = @(5,6)
  47 = @(6,6)
  48 = @(7,7)
  49 = @(5,5)
  50 = @(9,9)
  51 = @(7,7)
  52 = @(7,7)
  53 = @(10,10)
  54 = @(6,6)
  55 = @(7,7)
  56 = @(8,8)
  57 = @(6,6)
  58 = @(9,9)
  59 = @(8,8)
  60 = @(9,9)
  61 = @(6,7)
  62 = @(8,8)
  63 = @(7,7)
  64 = @(9,9)
  65 = @(6,7)
  66 = @(9,9)
  67 = @(8,8)
  68 = @(8,8)
  69 = @(8,8)
  70 = @(9,9)
  71 = @(10,10)
  72 = @(8,8)
  73 = @(6,6)
  74 = @(6,6)
  75 = @(4,4)
  76 = @(5,5)
}

foreach ($r in $data.Keys) {
  $vals = $data[$r]
  $ws.Cells.Item($r, 9).Value = $vals[0]
  $ws.Cells.Item($r, 10).Value = $vals[1]
}

Write-Host "Added columns I0 and IF"
